$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7 data (mirrors the pattern used in rows 2-6) ---
# Leading apostrophes force text/quote-prefixed entry so the cells pick up
# the same "quotePrefix" style (s="1") as their counterparts in rows 2-6,
# and so numeric-looking strings (phone numbers) stay text instead of
# becoming numbers.
$ws.Range("B7").Value = "'Chrome"
$ws.Range("C7").Value = "'9830252735"
$ws.Range("D7").Value = "'tekoda"
$ws.Range("E7").Value = "'BSNL Payment Portal"
$ws.Range("F7").Value = "Authorize | BSNL Portal"
$ws.Range("G7").Value = "'Dashboard Welcome to BSNL Payment Portal"

# --- Data validation: drop the disjoint "D7" piece, keep "B2:B7" ---
$ws.Range("D7").Validation.Delete()

# --- Column A width tweak ---
$ws.Columns.Item(1).ColumnWidth = 11.666666666666668
